$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.43159475247495
$ws.Range("C2").Value = 9.860254230117771
$ws.Range("E2").Value = 12.04184355042982
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 15.1515694203765
$ws.Range("H2").Value = 10.94917713737281
$ws.Range("I2").Value = 15.04975895049383
$ws.Range("M2").Value = 14.06048477013874
$ws.Range("O2").Value = 14.82721739443905
$ws.Range("B3").Value = 11.69621503460326
$ws.Range("C3").Value = 9.457720933776473
$ws.Range("E3").Value = 11.99570691534118
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 15.29675550950922
$ws.Range("H3").Value = 11.01490787016469
$ws.Range("I3").Value = 15.19632566768479
$ws.Range("M3").Value = 13.70630745484307
$ws.Range("O3").Value = 14.94586293071683
$ws.Range("B4").Value = 11.2184636832313
$ws.Range("C4").Value = 9.200749769805338
$ws.Range("E4").Value = 11.97327616963715
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 15.39707441954827
$ws.Range("H4").Value = 11.0577216290605
$ws.Range("I4").Value = 15.29080555170222
$ws.Range("M4").Value = 13.48611351455491
$ws.Range("O4").Value = 15.02385402747129
$ws.Range("B5").Value = 11.01724441895473
$ws.Range("C5").Value = 9.093680079066731
$ws.Range("E5").Value = 11.96562248142296
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 15.44072509831167
$ws.Range("H5").Value = 11.07578593190897
$ws.Range("I5").Value = 15.33043839927076
$ws.Range("M5").Value = 13.39582594924826
$ws.Range("O5").Value = 15.05692468296749
$ws.Range("B6").Value = 10.98343979282859
$ws.Range("C6").Value = 9.075762886453793
$ws.Range("E6").Value = 11.96444146712967
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 15.44813938297607
$ws.Range("H6").Value = 11.07882278451727
$ws.Range("I6").Value = 15.33708783775385
$ws.Range("M6").Value = 13.38080403754335
$ws.Range("O6").Value = 15.06249374505185
$ws.Range("B7").Value = 11.21577632099322
$ws.Range("C7").Value = 9.199315148477792
$ws.Range("E7").Value = 11.97316692520074
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 15.39765194561637
$ws.Range("H7").Value = 11.05796275084328
$ws.Range("I7").Value = 15.29133546845223
$ws.Range("M7").Value = 13.48489794390381
$ws.Range("O7").Value = 15.02429481794711
$ws.Range("B8").Value = 12.18351149233394
$ws.Range("C8").Value = 9.7235625572772
$ws.Range("E8").Value = 12.02471499368237
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 15.19928653530534
$ws.Range("H8").Value = 10.97133162857968
$ws.Range("I8").Value = 15.09936571382454
$ws.Range("M8").Value = 13.93900604126515
$ws.Range("O8").Value = 14.86705641298069
$ws.Range("B9").Value = 13.87116990619528
$ws.Range("C9").Value = 10.66938992328009
$ws.Range("E9").Value = 12.1722714704455
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 14.90071075633222
$ws.Range("H9").Value = 10.8209233404409
$ws.Range("I9").Value = 14.75837311558266
$ws.Range("M9").Value = 14.80253332012838
$ws.Range("O9").Value = 14.5997158163436
$ws.Range("B10").Value = 14.98123852134117
$ws.Range("C10").Value = 11.3092745727896
$ws.Range("E10").Value = 12.30837930136073
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 14.73881567089431
$ws.Range("H10").Value = 10.72228714064303
$ws.Range("I10").Value = 14.52925722134442
$ws.Range("M10").Value = 15.41389335977962
$ws.Range("O10").Value = 14.42858139807703
$ws.Range("B11").Value = 15.45784992116533
$ws.Range("C11").Value = 11.5876288793586
$ws.Range("E11").Value = 12.37611991026151
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 14.67812532409488
$ws.Range("H11").Value = 10.67999221986602
$ws.Range("I11").Value = 14.42963202536706
$ws.Range("M11").Value = 15.68576024761363
$ws.Range("O11").Value = 14.35628172681312
$ws.Range("B12").Value = 15.63424340035011
$ws.Range("C12").Value = 11.69114979433304
$ws.Range("E12").Value = 12.40258906137403
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 14.65704419797162
$ws.Range("H12").Value = 10.66434672773435
$ws.Range("I12").Value = 14.39256491020577
$ws.Range("M12").Value = 15.78771818603239
$ws.Range("O12").Value = 14.32970755275092
$ws.Range("B13").Value = 15.59643587017102
$ws.Range("C13").Value = 11.66893936969948
$ws.Range("E13").Value = 12.39685245653296
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 14.66149926456646
$ws.Range("H13").Value = 10.6676997720773
$ws.Range("I13").Value = 14.40051872130141
$ws.Range("M13").Value = 15.76580526277125
$ws.Range("O13").Value = 14.33539491429715
$ws.Range("B14").Value = 15.47244390627123
$ws.Range("C14").Value = 11.59618368627793
$ws.Range("E14").Value = 12.37828128886485
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 14.67635265786125
$ws.Range("H14").Value = 10.67869762462113
$ws.Range("I14").Value = 14.42656930439921
$ws.Range("M14").Value = 15.69416878458005
$ws.Range("O14").Value = 14.35407930215677
$ws.Range("B15").Value = 15.39596253993067
$ws.Range("C15").Value = 11.55137165613582
$ws.Range("E15").Value = 12.36701170039084
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 14.68569944912267
$ws.Range("H15").Value = 10.6854824116399
$ws.Range("I15").Value = 14.44261174935691
$ws.Range("M15").Value = 15.65015750324763
$ws.Range("O15").Value = 14.36562893418926
$ws.Range("B16").Value = 14.94951909000745
$ws.Range("C16").Value = 11.29082236305119
$ws.Range("E16").Value = 12.30406781548412
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 14.74304592136729
$ws.Range("H16").Value = 10.72510306281298
$ws.Range("I16").Value = 14.53586030129428
$ws.Range("M16").Value = 15.39599233426824
$ws.Range("O16").Value = 14.43341854749186
$ws.Range("B17").Value = 14.66836671367588
$ws.Range("C17").Value = 11.12768086969177
$ws.Range("E17").Value = 12.2669318860418
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 14.78157206914276
$ws.Range("H17").Value = 10.75006886163197
$ws.Range("I17").Value = 14.59424163225671
$ws.Range("M17").Value = 15.23840140601701
$ws.Range("O17").Value = 14.47643105785886
$ws.Range("B18").Value = 14.50398573264281
$ws.Range("C18").Value = 11.03265045250665
$ws.Range("E18").Value = 12.24612188165492
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 14.80494865404842
$ws.Range("H18").Value = 10.76467089196663
$ws.Range("I18").Value = 14.6282542415677
$ws.Range("M18").Value = 15.14717669679027
$ws.Range("O18").Value = 14.50169275884776
$ws.Range("B19").Value = 14.44787105477668
$ws.Range("C19").Value = 11.00027123193332
$ws.Range("E19").Value = 12.23917094550329
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 14.81307146646018
$ws.Range("H19").Value = 10.76965650422495
$ws.Range("I19").Value = 14.63984481715131
$ws.Range("M19").Value = 15.11619252379663
$ws.Range("O19").Value = 14.51033541745032
$ws.Range("B20").Value = 14.69857224501044
$ws.Range("C20").Value = 11.14517174038799
$ws.Range("E20").Value = 12.27082832711812
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 14.77734462270626
$ws.Range("H20").Value = 10.74738612352265
$ws.Range("I20").Value = 14.58798202448724
$ws.Range("M20").Value = 15.2552382320547
$ws.Range("O20").Value = 14.47179823015346
$ws.Range("B21").Value = 15.50897437847497
$ws.Range("C21").Value = 11.61760537018612
$ws.Range("E21").Value = 12.38371408432986
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 14.67193797083735
$ws.Range("H21").Value = 10.67545722550663
$ws.Range("I21").Value = 14.41889975960416
$ws.Range("M21").Value = 15.71523779344218
$ws.Range("O21").Value = 14.34856937142007
$ws.Range("B22").Value = 16.01478449345331
$ws.Range("C22").Value = 11.91535637540888
$ws.Range("E22").Value = 12.46224400767994
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 14.61414532309532
$ws.Range("H22").Value = 10.63060832429414
$ws.Range("I22").Value = 14.31223334917795
$ws.Range("M22").Value = 16.01005291488666
$ws.Range("O22").Value = 14.27272226918735
$ws.Range("B23").Value = 15.74700701664974
$ws.Range("C23").Value = 11.75746461908303
$ws.Range("E23").Value = 12.41990354992348
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 14.64396295880759
$ws.Range("H23").Value = 10.65434717995178
$ws.Range("I23").Value = 14.36881290226884
$ws.Range("M23").Value = 15.85326627897333
$ws.Range("O23").Value = 14.31277211496909
$ws.Range("B24").Value = 14.68492486652292
$ws.Range("C24").Value = 11.1372679702439
$ws.Range("E24").Value = 12.26906506274346
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 14.77925203091153
$ws.Range("H24").Value = 10.74859821425157
$ws.Range("I24").Value = 14.59081059616832
$ws.Range("M24").Value = 15.24762823808107
$ws.Range("O24").Value = 14.47389107099682
$ws.Range("B25").Value = 13.43732677375843
$ws.Range("C25").Value = 10.4228731318352
$ws.Range("E25").Value = 12.12743262474492
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 14.97155779398275
$ws.Range("H25").Value = 10.85952838282768
$ws.Range("I25").Value = 14.84684584357032
$ws.Range("M25").Value = 14.57252858849783
$ws.Range("O25").Value = 14.66762014845584
